# Actualiza base de datos EC: renumera los periodos de mora (columna E) en
# orden ascendente (1610 -> 2003) y actualiza los valores de Valor Mora
# (columna F) y Salario Basico (columna G) para las filas 16 a 57.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Hoja1")
if (-not $ws) {
    $ws = $wb.ActiveSheet
}

# Nuevas etiquetas de periodo, en orden ascendente, para las filas 16..57
$periods = @(
    "1610","1611","1612",
    "1701","1702","1703","1704","1705","1706","1707","1708","1709","1710","1711","1712",
    "1801","1802","1803","1804","1805","1806","1807","1808","1809","1810","1811","1812",
    "1901","1902","1903","1904","1905","1906","1907","1908","1909","1910","1911","1912",
    "2001","2002","2003"
)

$startRow = 16
for ($i = 0; $i -lt $periods.Length; $i++) {
    $row = $startRow + $i
    $ws.Cells.Item($row, 5).Value = $periods[$i]

    if ($row -le 38) {
        $ws.Cells.Item($row, 6).Value = 27578
    } else {
        $ws.Cells.Item($row, 6).Value = 31249
    }

    $ws.Cells.Item($row, 7).Value = 781242
}
